{"js": "const body = context.document.body;\n\n// ---------------------------------------------------------------------------\n// 1) Title paragraph: \"Restaurante \"Conquistador Caf\u00e9\"\" -> \"Restaurante \"Caf\u00e9 Conquistador\"\"\n// ---------------------------------------------------------------------------\nconst titleHits = body.search('Conquistador Caf\u00e9', { matchCase: true });\ntitleHits.load('items');\nawait context.sync();\n\nif (titleHits.items.length > 0) {\n  titleHits.items[0].insertText('Caf\u00e9 Conquistador', 'Replace');\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------------\n// 2) Body paragraph: insert the ID clause right before \"se compromete a utilizar\"\n// ---------------------------------------------------------------------------\nconst clauseHits = body.search('se compromete a utilizar y probar el sistema', { matchCase: true });\nclauseHits.load('items');\nawait context.sync();\n\nif (clauseHits.items.length > 0) {\n  clauseHits.items[0].insertText(\n    'con n\u00famero de cedula de ciudadan\u00eda 050016186-4 ',\n    'Before'\n  );\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------------\n// 3) Move the \"_GoBack\" bookmark from the signature underline paragraph to the\n//    end of the title paragraph's text (Word keeps only one \"_GoBack\" bookmark\n//    alive at a time - it tracks the most recent edit position).\n// ---------------------------------------------------------------------------\ncontext.document.deleteBookmark('_GoBack');\nawait context.sync();\n\nconst titleAgainHits = body.search('Restaurante \"Caf\u00e9 Conquistador\"', { matchCase: true });\ntitleAgainHits.load('items');\nawait context.sync();\n\nif (titleAgainHits.items.length > 0) {\n  const titlePara = titleAgainHits.items[0].paragraphs.getFirst();\n  const titleContent = titlePara.getRange('Content');\n  titleContent.insertBookmark('_GoBack');\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------------\n// 4) Add a new centered paragraph \"Telf.: 0992963818\" right after \"Firma\".\n// ---------------------------------------------------------------------------\nconst firmaHits = body.search('Firma', { matchCase: true });\nfirmaHits.load('items');\nawait context.sync();\n\nif (firmaHits.items.length > 0) {\n  const firmaPara = firmaHits.items[0].paragraphs.getFirst();\n  const phonePara = firmaPara.insertParagraph('Telf.: 0992963818', 'After');\n  phonePara.alignment = 'Centered';\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# 1) Title paragraph: \"Restaurante \"Conquistador Caf\u00e9\"\" -> \"Restaurante \"Caf\u00e9 Conquistador\"\"\n# ---------------------------------------------------------------------------\n$titleRange = $d.Content\n$titleRange.Find.ClearFormatting()\n$titleFound = $titleRange.Find.Execute('Conquistador Caf\u00e9')\nif ($titleFound) {\n    $titleRange.Text = \"Caf\u00e9 Conquistador\"\n}\n\n# ---------------------------------------------------------------------------\n# 2) Body paragraph: insert the ID clause right before \"se compromete a utilizar\"\n# ---------------------------------------------------------------------------\n$insertRange = $d.Content\n$insertRange.Find.ClearFormatting()\n$insertFound = $insertRange.Find.Execute('se compromete a utilizar y probar el sistema')\nif ($insertFound) {\n    $insertRange.Collapse(1)\n    $insertRange.InsertBefore(\"con n\u00famero de cedula de ciudadan\u00eda 050016186-4 \")\n}\n\n# ---------------------------------------------------------------------------\n# 3) Move the \"_GoBack\" bookmark from the signature underline paragraph to the\n#    end of the title paragraph's text (Word keeps only one \"_GoBack\" bookmark\n#    alive at a time - it tracks the most recent edit position).\n# ---------------------------------------------------------------------------\ntry {\n    $oldBm = $d.Bookmarks(\"_GoBack\")\n    $oldBm.Delete()\n} catch {\n}\n\n$titleAgain = $d.Content\n$titleAgain.Find.ClearFormatting()\n$titleFoundAgain = $titleAgain.Find.Execute('Restaurante \"Caf\u00e9 Conquistador\"')\nif ($titleFoundAgain) {\n    $titlePara = $titleAgain.Paragraphs(1)\n    $bmRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n\n# ---------------------------------------------------------------------------\n# 4) Add a new centered paragraph \"Telf.: 0992963818\" right after \"Firma\".\n# ---------------------------------------------------------------------------\n$firmaRange = $d.Content\n$firmaRange.Find.ClearFormatting()\n$firmaFound = $firmaRange.Find.Execute(\"Firma\")\nif ($firmaFound) {\n    $firmaPara = $firmaRange.Paragraphs(1)\n    $firmaPara.Range.InsertParagraphAfter()\n    $phonePara = $firmaPara.Next()\n    $phonePara.Range.Text = \"Telf.: 0992963818\"\n    $phonePara.Range.ParagraphFormat.Alignment = 1\n}\n"}
